# Fix parser for notes: the footnote markers (*, **, ***, ****) stored as
# shared-string text in column L are replaced with their corresponding
# plain numeric values (1, 2, 3, 4) so the notes parser can read them as
# numbers instead of literal asterisk strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

$ws.Range("L6").Value  = 1
$ws.Range("L7").Value  = 1
$ws.Range("L8").Value  = 2
$ws.Range("L9").Value  = 2
$ws.Range("L10").Value = 3
$ws.Range("L11").Value = 4
$ws.Range("L12").Value = 4
$ws.Range("L15").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("L18").Value = 1

# Restore the on-screen selection the author left behind when saving.
$ws.Range("L13").Select() | Out-Null
